$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- header text updates (Volume number, date range) ----
$ws.Range("A8").Value = "Volume 31   Number  4"
$ws.Range("C9").Value = "Report Covering the Week  1/22/2024  Through  1/28/2024"

# ---- simple numeric value updates ----
$ws.Range("G14").Value = 1
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = -66.666666666666
$ws.Range("J15").Value = 3
$ws.Range("K15").Value = -66.666666666666
$ws.Range("C16").Value = 8
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = 14.285714285714
$ws.Range("F16").Value = 28
$ws.Range("G16").Value = 25
$ws.Range("H16").Value = 12
$ws.Range("I16").Value = 28
$ws.Range("J16").Value = 25
$ws.Range("K16").Value = 12
$ws.Range("L16").Value = 115.384615384615
$ws.Range("M16").Value = -17.647058823529
$ws.Range("N16").Value = -65.432098765432
$ws.Range("D17").Value = 12
$ws.Range("E17").Value = -33.333333333333
$ws.Range("G17").Value = 25
$ws.Range("H17").Value = 28
$ws.Range("I17").Value = 32
$ws.Range("J17").Value = 25
$ws.Range("K17").Value = 28
$ws.Range("L17").Value = 45.454545454545
$ws.Range("M17").Value = 68.421052631578
$ws.Range("N17").Value = 146.153846153846
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 4
$ws.Range("F18").Value = 14
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = -30
$ws.Range("I18").Value = 14
$ws.Range("J18").Value = 20
$ws.Range("K18").Value = -30
$ws.Range("L18").Value = -6.666666666666
$ws.Range("M18").Value = -44
$ws.Range("N18").Value = -90.410958904109
$ws.Range("C19").Value = 25
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 257.142857142857
$ws.Range("F19").Value = 73
$ws.Range("H19").Value = 55.31914893617
$ws.Range("I19").Value = 73
$ws.Range("J19").Value = 47
$ws.Range("K19").Value = 55.31914893617
$ws.Range("L19").Value = 10.60606060606
$ws.Range("M19").Value = 143.333333333333
$ws.Range("N19").Value = 82.5
$ws.Range("C20").Value = 14
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 39
$ws.Range("G20").Value = 36
$ws.Range("H20").Value = 8.333333333333
$ws.Range("I20").Value = 39
$ws.Range("J20").Value = 36
$ws.Range("K20").Value = 8.333333333333
$ws.Range("L20").Value = 44.444444444444
$ws.Range("M20").Value = 200
$ws.Range("N20").Value = -71.739130434782
$ws.Range("C21").Value = 57
$ws.Range("D21").Value = 39
$ws.Range("E21").Value = 46.153846153846
$ws.Range("F21").Value = 187
$ws.Range("G21").Value = 157
$ws.Range("H21").Value = 19.108280254777
$ws.Range("I21").Value = 187
$ws.Range("J21").Value = 157
$ws.Range("K21").Value = 19.108280254777
$ws.Range("L21").Value = 28.082191780821
$ws.Range("M21").Value = 52.032520325203
$ws.Range("N21").Value = -55.581947743467
$ws.Range("E22").Value = -100
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -33.333333333333
$ws.Range("J22").Value = 3
$ws.Range("K22").Value = -33.333333333333
$ws.Range("L22").Value = 0
$ws.Range("C23").Value = 4
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 100
$ws.Range("F23").Value = 11
$ws.Range("G23").Value = 13
$ws.Range("H23").Value = -15.384615384615
$ws.Range("I23").Value = 11
$ws.Range("J23").Value = 13
$ws.Range("K23").Value = -15.384615384615
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 57.142857142857
$ws.Range("C24").Value = 38
$ws.Range("D24").Value = 50
$ws.Range("E24").Value = -24
$ws.Range("F24").Value = 98
$ws.Range("G24").Value = 132
$ws.Range("H24").Value = -25.757575757575
$ws.Range("I24").Value = 98
$ws.Range("J24").Value = 132
$ws.Range("K24").Value = -25.757575757575
$ws.Range("L24").Value = 12.643678160919
$ws.Range("M24").Value = 10.112359550561
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = -20
$ws.Range("F25").Value = 34
$ws.Range("G25").Value = 38
$ws.Range("H25").Value = -10.526315789473
$ws.Range("I25").Value = 34
$ws.Range("J25").Value = 38
$ws.Range("K25").Value = -10.526315789473
$ws.Range("L25").Value = -12.820512820512
$ws.Range("M25").Value = -10.526315789473
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = -100
$ws.Range("F26").Value = 2
$ws.Range("G26").Value = 5
$ws.Range("H26").Value = -60
$ws.Range("J26").Value = 5
$ws.Range("K26").Value = -60
$ws.Range("F27").Value = 4
$ws.Range("H27").Value = -42.857142857142
$ws.Range("I27").Value = 4
$ws.Range("K27").Value = -42.857142857142
$ws.Range("L27").Value = 33.333333333333
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = -50
$ws.Range("L28").Value = -83.333333333333
$ws.Range("N28").Value = -50
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = -50
$ws.Range("L29").Value = -80
$ws.Range("N29").Value = -50

# ---- text -> number conversions (set value + numberformat) ----
$ws.Range("L14").Value = -100
$ws.Range("L14").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("D15").Value = 2
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("E15").Value = -100
$ws.Range("E15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("C27").Value = 1
$ws.Range("C27").NumberFormat = "#,##0"

# ---- number -> text conversions (copy format+value from stable donor cells) ----
$ws.Range("C30").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C30").Copy()
$ws.Range("C22").PasteSpecial(-4163)
$ws.Range("C30").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("C30").Copy()
$ws.Range("C26").PasteSpecial(-4163)
$ws.Range("C30").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("C30").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E30").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("E30").Copy()
$ws.Range("E27").PasteSpecial(-4163)
$ws.Range("C30").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C30").Copy()
$ws.Range("C28").PasteSpecial(-4163)
$ws.Range("C30").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("C30").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E30").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E30").Copy()
$ws.Range("E28").PasteSpecial(-4163)
$ws.Range("C30").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("C30").Copy()
$ws.Range("C29").PasteSpecial(-4163)
$ws.Range("C30").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("C30").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E30").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("E30").Copy()
$ws.Range("E29").PasteSpecial(-4163)

Write-Output "edit complete"
